{"js": "// Resume edit: flesh out the \"Environmental Specialist II\" bullet points\n// (City of Los Angeles, Public Works Office of Forest Management entry).\n\nconst body = context.document.body;\n\n// Helper: find the first exact match for `needle` under `scope` (a Body or\n// Range) and replace it with `replacement`, preserving the formatting of\n// the run(s) the match lands on.\nasync function replaceOnce(scope, needle, replacement) {\n  const results = scope.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + needle);\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Develop the City's first Urban Forest Management Plan\" -> add a\n//    trailing period.\nawait replaceOnce(\n  body,\n  \"Develop the City\\u2019s first Urban Forest Management Plan\",\n  \"Develop the City\\u2019s first Urban Forest Management Plan.\"\n);\n\n// 2) \"...Office of Forest Management (OFM) to \" -> drop the dangling \"to\"\n//    and add the rest of the sentence.\nawait replaceOnce(\n  body,\n  \"Provide technical expertise related to urban forestry and data science for the Office of Forest Management (OFM) to \",\n  \"Provide technical expertise related to urban forestry and data science for the Office of Forest Management (OFM), other City departments, and to external City partners on City policies and projects. \"\n);\n\n// 3) \"Supported urban forestry research efforts within the Los Angeles\n//    region \" -> change tense and append a new clause.\nawait replaceOnce(\n  body,\n  \"Supported urban forestry research efforts within the Los Angeles region \",\n  \"Support urban forestry research efforts within the Los Angeles region, with an emphasis on canopy equity.\"\n);\n", "ps1": "# Resume edit: flesh out the \"Environmental Specialist II\" bullet points\n# (City of Los Angeles, Public Works Office of Forest Management entry).\n\n$d = $word.ActiveDocument\n\n# NOTE: named parameters (e.g. \"-FindText foo\") are not reliable in this\n# PowerShell host, so the helper below takes positional arguments.\nfunction Replace-ExactText($FindText, $ReplaceText) {\n    $range = $d.Content\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $found = $range.Find.Execute($FindText)\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n    $range.Text = $ReplaceText\n}\n\n$rsquo = [char]0x2019\n\n# 1) \"Develop the City's first Urban Forest Management Plan\" -> add a\n#    trailing period.\nReplace-ExactText \"Develop the City${rsquo}s first Urban Forest Management Plan\" \"Develop the City${rsquo}s first Urban Forest Management Plan.\"\n\n# 2) \"...Office of Forest Management (OFM) to \" -> drop the dangling \"to\"\n#    and add the rest of the sentence.\nReplace-ExactText \"Provide technical expertise related to urban forestry and data science for the Office of Forest Management (OFM) to \" \"Provide technical expertise related to urban forestry and data science for the Office of Forest Management (OFM), other City departments, and to external City partners on City policies and projects. \"\n\n# 3) \"Supported urban forestry research efforts within the Los Angeles\n#    region \" -> change tense and append a new clause.\nReplace-ExactText \"Supported urban forestry research efforts within the Los Angeles region \" \"Support urban forestry research efforts within the Los Angeles region, with an emphasis on canopy equity.\"\n"}
